# Generate Report for Handback
# Update the "Latest HO Xliff Generate Date" / "Correspond Handoff Datetime" /
# "Correspond Handback DateTime" values for the file
# "a6be5ba1-4d7a-457f-9afe-923b2a13a896.md" (row 4 on every sheet) to reflect a
# fresh handback report generation.

$wb = $excel.ActiveWorkbook

# --- Overview sheet ---
$wsOverview = $wb.Worksheets.Item("Overview")
$wsOverview.Range("G4").Value = "2017-02-21 04:15:31"

# --- zh-cn sheet ---
$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsZhCn.Range("H4").Value = "2017-02-21 04:15:15"
$wsZhCn.Range("L4").Value = "2017-02-21 04:16:10"

# --- de-de sheet ---
$wsDeDe = $wb.Worksheets.Item("de-de")
$wsDeDe.Range("H4").Value = "2017-02-21 04:15:31"
$wsDeDe.Range("L4").Value = "2017-02-21 04:16:33"
